$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.78
$ws.Range("G3").Value = 1.46
$ws.Range("H3").Value = 8.4
$ws.Range("K3").Value = 5.3
$ws.Range("L3").Value = 1.31
$ws.Range("P3").Value = 2.48
$ws.Range("Q3").Value = 1.65
$ws.Range("R3").Value = 1.58
$ws.Range("S3").Value = 2.64
$ws.Range("T3").Value = 1.9
$ws.Range("U3").Value = 2.04
$ws.Range("W3").Value = 3.2
$ws.Range("N4").Value = 5.2
$ws.Range("R4").Value = 1.56
$ws.Range("S4").Value = 2.68
$ws.Range("F5").Value = 1.5
$ws.Range("G5").Value = 1.67
$ws.Range("H5").Value = 4.7
$ws.Range("I5").Value = 6.6
$ws.Range("K5").Value = 5.9
$ws.Range("O5").Value = 1.13
$ws.Range("P5").Value = 2.74
$ws.Range("Q5").Value = 1.4
$ws.Range("R5").Value = 1.73
$ws.Range("S5").Value = 1.95
$ws.Range("T5").Value = 1.53
$ws.Range("U5").Value = 2.44
$ws.Range("W5").Value = 2.48
$ws.Range("X5").Value = 40
$ws.Range("Z5").Value = 1000
$ws.Range("AB5").Value = 17
$ws.Range("AC5").Value = 14
$ws.Range("AE5").Value = 65
$ws.Range("AF5").Value = 15.5
$ws.Range("AI5").Value = 55
$ws.Range("AN5").Value = 5.5
$ws.Range("H6").Value = 2.4
$ws.Range("Q7").Value = 1.71
$ws.Range("R7").Value = 1.45
$ws.Range("S7").Value = 2.8
$ws.Range("L8").Value = 1.46
$ws.Range("O8").Value = 1.37
$ws.Range("P8").Value = 1.83
$ws.Range("AN8").Value = 24
$ws.Range("F9").Value = 2.74
$ws.Range("G9").Value = 2.8
$ws.Range("H9").Value = 2.54
$ws.Range("I9").Value = 2.56
$ws.Range("L9").Value = 1.27
$ws.Range("O9").Value = 1.16
$ws.Range("P9").Value = 2.82
$ws.Range("Q9").Value = 1.52
$ws.Range("R9").Value = 1.74
$ws.Range("S9").Value = 2.28
$ws.Range("U9").Value = 2.96
$ws.Range("V9").Value = 1.64
$ws.Range("W9").Value = 1.54
$ws.Range("X9").Value = 26
$ws.Range("Y9").Value = 18
$ws.Range("AB9").Value = 20
$ws.Range("AC9").Value = 9.800000000000001
$ws.Range("AE9").Value = 22
$ws.Range("AF9").Value = 24
$ws.Range("AI9").Value = 26
$ws.Range("AJ9").Value = 44
$ws.Range("AK9").Value = 24
$ws.Range("AL9").Value = 28
$ws.Range("AM9").Value = 50
$ws.Range("AN9").Value = 13.5
$ws.Range("AO9").Value = 12
$ws.Range("F10").Value = 9.6
$ws.Range("G10").Value = 10
$ws.Range("J10").Value = 5.7
$ws.Range("L10").Value = 1.3
$ws.Range("N10").Value = 5.5
$ws.Range("P10").Value = 2.52
$ws.Range("T10").Value = 1.96
$ws.Range("W10").Value = 1.11
$ws.Range("AJ10").Value = 320
$ws.Range("AN10").Value = 140
$ws.Range("I11").Value = 24
$ws.Range("N11").Value = 6.8
$ws.Range("AF11").Value = 7.8
$ws.Range("AH11").Value = 150
$ws.Range("AJ11").Value = 8.4
$ws.Range("F12").Value = 1.31
$ws.Range("G12").Value = 1.32
$ws.Range("H12").Value = 10.5
$ws.Range("L12").Value = 1.22
$ws.Range("N12").Value = 8.4
$ws.Range("U12").Value = 2.3
$ws.Range("W12").Value = 4.1
$ws.Range("AA12").Value = 360
$ws.Range("AH12").Value = 24
$ws.Range("AM12").Value = 90
$ws.Range("F13").Value = 5.7
$ws.Range("G13").Value = 5.9
$ws.Range("H13").Value = 1.68
$ws.Range("I13").Value = 1.7
$ws.Range("J13").Value = 4.2
$ws.Range("K13").Value = 4.3
$ws.Range("L13").Value = 1.37
$ws.Range("N13").Value = 4.4
$ws.Range("R13").Value = 1.44
$ws.Range("S13").Value = 3.15
$ws.Range("T13").Value = 1.88
$ws.Range("U13").Value = 2.08
$ws.Range("V13").Value = 2.42
$ws.Range("W13").Value = 1.2
$ws.Range("X13").Value = 17
$ws.Range("Y13").Value = 9.4
$ws.Range("Z13").Value = 9.800000000000001
$ws.Range("AA13").Value = 16.5
$ws.Range("AB13").Value = 21
$ws.Range("AC13").Value = 9.199999999999999
$ws.Range("AD13").Value = 9.800000000000001
$ws.Range("AF13").Value = 48
$ws.Range("AG13").Value = 22
$ws.Range("AH13").Value = 20
$ws.Range("AJ13").Value = 160
$ws.Range("AK13").Value = 75
$ws.Range("AL13").Value = 75
$ws.Range("AN13").Value = 85
$ws.Range("AO13").Value = 8.800000000000001
$ws.Range("G14").Value = 3.2
$ws.Range("I14").Value = 2.44
$ws.Range("N14").Value = 4.9
$ws.Range("T14").Value = 1.63
$ws.Range("V14").Value = 1.69
$ws.Range("AJ14").Value = 50
$ws.Range("AO14").Value = 14.5
$ws.Range("F15").Value = 2.28
$ws.Range("G15").Value = 2.9
$ws.Range("H15").Value = 2.48
$ws.Range("I15").Value = 3.1
$ws.Range("K15").Value = 4.6
$ws.Range("L15").Value = 1.25
$ws.Range("M15").Value = 1.04
$ws.Range("V15").Value = 1.48
$ws.Range("W15").Value = 1.52
$ws.Range("Z15").Value = 25
$ws.Range("AD15").Value = 15
$ws.Range("AE15").Value = 32
$ws.Range("AI15").Value = 38
$ws.Range("AN15").Value = 17
$ws.Range("AO15").Value = 21
$ws.Range("F16").Value = 2.5
$ws.Range("G16").Value = 2.86
$ws.Range("I16").Value = 3.25
$ws.Range("N16").Value = 3.35
$ws.Range("S16").Value = 3.6
$ws.Range("T16").Value = 1.76
$ws.Range("V16").Value = 1.45
$ws.Range("W16").Value = 1.58
$ws.Range("Y16").Value = 14
$ws.Range("AC16").Value = 9.199999999999999
$ws.Range("AF16").Value = 21
$ws.Range("AJ16").Value = 46
